$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2ND Q 2023")
$ws.Rows.Item(14).Copy($ws.Rows.Item(15))
Write-Output $ws.Range("A15").Value()
Write-Output $ws.Range("H15").Value()
